# Apply "F" column (想去人数 / interested count) updates across sheets,
# matching the regenerated gh-pages data snapshot.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions) sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 7956
$ws1.Range("F14").Value = 38
$ws1.Range("F15").Value = 684
$ws1.Range("F22").Value = 11222
$ws1.Range("F24").Value = 87
$ws1.Range("F25").Value = 2136
$ws1.Range("F26").Value = 2906
$ws1.Range("F29").Value = 2543
$ws1.Range("F35").Value = 332
$ws1.Range("F38").Value = 70
$ws1.Range("F39").Value = 5665
$ws1.Range("F40").Value = 73
$ws1.Range("F42").Value = 806
$ws1.Range("F47").Value = 1481

# 本地生活 (Local life) sheet
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 166
$ws3.Range("F3").Value = 280

# 全部类型 (All types, combined) sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 166
$ws4.Range("F5").Value  = 280
$ws4.Range("F8").Value  = 7956
$ws4.Range("F17").Value = 684
$ws4.Range("F25").Value = 11222
$ws4.Range("F27").Value = 87
$ws4.Range("F28").Value = 2136
$ws4.Range("F29").Value = 2906
$ws4.Range("F30").Value = 2543
$ws4.Range("F36").Value = 332
$ws4.Range("F39").Value = 70
$ws4.Range("F40").Value = 5665
$ws4.Range("F43").Value = 806
$ws4.Range("F48").Value = 1481
